# feat: add 2022-Q4 data
#
# 1) Insert a brand-new worksheet named "2022-Q4" right after "总计",
#    carrying the newest per-fund holdings table (the data that used to
#    live on the "2022-Q3" tab gets superseded there and moves, unchanged,
#    one tab to the right - handled automatically since we only ever
#    insert a sheet and edit "总计"; every other existing sheet/tab keeps
#    its own name+content and simply shifts position).
# 2) Update the "总计" (totals) sheet: add a "2022-Q4" row at the top of
#    the data and push the existing quarters down by one row, appending
#    the oldest quarter ("2021-Q1") as the new last row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet - shift all quarter rows down by one and
# insert the new 2022-Q4 totals at the top.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Make room for the new trailing row (2021-Q1) by copying the format of
# the current last data row (2021-Q1, row 7) down into the new row 8.
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)
$summary.Range("A8").Value = 6

$quarters = @(
    @("2022-Q4", 4, 0.01),
    @("2022-Q3", 4, 0.12),
    @("2022-Q2", 1, 0.03),
    @("2022-Q1", 3, 0.04),
    @("2021-Q4", 4, 0.08),
    @("2021-Q3", 2, 0.04),
    @("2021-Q1", 4, 0.18)
)

for ($i = 0; $i -lt $quarters.Length; $i++) {
    $row = $i + 2
    $summary.Cells.Item($row, 2).Value = $quarters[$i][0]
    $summary.Cells.Item($row, 3).Value = $quarters[$i][1]
    $summary.Cells.Item($row, 4).Value = $quarters[$i][2]
}

# ---------------------------------------------------------------------
# Part 2: brand-new "2022-Q4" worksheet with the fund-level breakdown,
# placed right after "总计" (i.e. as the new second tab).
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $afterSheet)
$q4.Name = "2022-Q4"

# Use the (still unshifted) "2022-Q3" sheet purely as a formatting
# template - copy its header style and its column-A index style, then
# overwrite with the new text/values.
$template = $wb.Worksheets.Item("2022-Q3")

$template.Range("B1:H1").Copy()
$q4.Range("B1").PasteSpecial(-4122)

$template.Range("A2:A5").Copy()
$q4.Range("A2").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1
$q4.Range("A4").Value = 2
$q4.Range("A5").Value = 3

# Columns B..G hold text in the source workbook (fund codes/names and
# numeric-looking figures kept as strings) - force text formatting
# before assigning so values like "014214" keep their leading zero and
# "0.32" doesn't get silently coerced into a float.
$q4.Range("B2:G5").NumberFormat = "@"

$q4.Range("B2").Value = "014214"
$q4.Range("C2").Value = "光大保德信核心资产混合A"
$q4.Range("D2").Value = "0.32"
$q4.Range("E2").Value = "85.79"
$q4.Range("F2").Value = "2.87"
$q4.Range("G2").Value = "0.0092"
$q4.Range("H2").Value = 10

$q4.Range("B3").Value = "013182"
$q4.Range("C3").Value = "安信港股通精选混合C"
$q4.Range("D3").Value = "0.12"
$q4.Range("E3").Value = "69.28"
$q4.Range("F3").Value = "2.26"
$q4.Range("G3").Value = "0.0027"
$q4.Range("H3").Value = 10

$q4.Range("B4").Value = "014215"
$q4.Range("C4").Value = "光大保德信核心资产混合C"
$q4.Range("D4").Value = "0.03"
$q4.Range("E4").Value = "85.79"
$q4.Range("F4").Value = "2.87"
$q4.Range("G4").Value = "0.0009"
$q4.Range("H4").Value = 10

$q4.Range("B5").Value = "013181"
$q4.Range("C5").Value = "安信港股通精选混合A"
$q4.Range("D5").Value = "0.02"
$q4.Range("E5").Value = "69.28"
$q4.Range("F5").Value = "2.26"
$q4.Range("G5").Value = "0.0005"
$q4.Range("H5").Value = 10

# Drop the text number-format again so the cells fall back to the
# workbook's default style (matching the plain, un-styled data cells
# used throughout the rest of the workbook) instead of keeping an
# explicit "@" text format applied.
$q4.Range("B2:G5").ClearFormats()
